# Script 1 - atualização em 2025-09-20 17:07:42Z
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 12.73287554584524
$ws.Range("D5").Value = 13.10014696707147
$ws.Range("D6").Value = 11.94223451277554
$ws.Range("D7").Value = 11.76881496409633
$ws.Range("D8").Value = 12.10466695575146
$ws.Range("D9").Value = 12.74924934636568
$ws.Range("D10").Value = 12.09431389892419
$ws.Range("D11").Value = 13.82986703080613
$ws.Range("D12").Value = 8.675532766104695
$ws.Range("D13").Value = 16.73480696636505
$ws.Range("D16").Value = 14.72565595323923
$ws.Range("D17").Value = 12.63416987646953
$ws.Range("D18").Value = 12.4691647775766
$ws.Range("D19").Value = 12.47276728182001
$ws.Range("D20").Value = 12.57007644680532
$ws.Range("D22").Value = 13.52047189388596
$ws.Range("D23").Value = 9.583461848945253
$ws.Range("D34").Value = 4.560678164550927
$ws.Range("E34").Value = 23
